$wb = $excel.ActiveWorkbook

# "Greece" is a new test-data sheet for the Greece Market, cloned from the
# existing "Croatia" sheet template (same layout/styles), with the two
# market-specific cells (Jira ref + market name) updated.
$croatia = $wb.Worksheets.Item("Croatia")
$croatia.Copy($null, $croatia) | Out-Null

$greece = $wb.Worksheets.Item("Croatia (2)")
$greece.Name = "Greece"

# Write B4 before B2 so new shared-string entries land in the same order
# as the authored workbook (ticket ref first, market name second).
$greece.Range("B4").Value = "NGC-4119/T3169"
$greece.Range("B2").Value = "Greece Market"

# Restore cursor positions / active sheet like the authored workbook.
$croatia.Range("G28").Select() | Out-Null
$greece.Range("H27").Select() | Out-Null
$greece.Activate() | Out-Null
